# Update the quarterly (yearly period) balance-sheet database: roll the
# 5-period window forward by one period (drop the oldest period, shift the
# remaining four left, and append the newest period's figures + labels).
# Mirrors the author commit "update database and change read_price algorithm".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Period / publish-date header labels (row 8 = دوره مالی, row 9 = تاریخ انتشار) ---
$ws.Range('D8').Value = '12 ماهه منتهی به 1397/12'
$ws.Range('E8').Value = '12 ماهه منتهی به 1398/12'
$ws.Range('F8').Value = '12 ماهه منتهی به 1399/12'
$ws.Range('G8').Value = '12 ماهه منتهی به 1400/12'
$ws.Range('H8').Value = '12 ماهه منتهی به 1401/12'
$ws.Range('D9').Value = '1399-02-31 (10)'
$ws.Range('E9').Value = '1400-02-30 (8)'
$ws.Range('F9').Value = '1401-02-31 (8)'
$ws.Range('G9').Value = '1402-02-10 (7)'
# H9 ("1402-02-10") parses as a date literal via normal assignment, so force
# text format first, then restore the original (General) cell style by
# pasting formats from the never-edited C9 template cell (same style as H9).
$ws.Range('H9').NumberFormat = '@'
$ws.Range('H9').Value = '1402-02-10'
$ws.Range('C9').Copy() | Out-Null
$ws.Range('H9').PasteSpecial(-4122) | Out-Null
$ws.Range('D12').Value = 2736538
$ws.Range('E12').Value = 5842571
$ws.Range('F12').Value = 3688510
$ws.Range('G12').Value = 3515523
$ws.Range('H12').Value = 2370001
$ws.Range('D13').Value = 1865000
$ws.Range('E13').Value = 2065000
$ws.Range('F13').Value = 6234030
$ws.Range('G13').Value = 14384068
$ws.Range('H13').Value = 23222503
$ws.Range('D14').Value = 1577746
$ws.Range('E14').Value = 2516928
$ws.Range('F14').Value = 8544737
$ws.Range('G14').Value = 20611130
$ws.Range('H14').Value = 32045778
$ws.Range('D15').Value = 1111555
$ws.Range('E15').Value = 1025403
$ws.Range('F15').Value = 1085998
$ws.Range('G15').Value = 3699041
$ws.Range('H15').Value = 2250423
$ws.Range('D16').Value = 64539
$ws.Range('E16').Value = 59382
$ws.Range('F16').Value = 199065
$ws.Range('G16').Value = 283888
$ws.Range('H16').Value = 472859
$ws.Range('D18').Value = 7355378
$ws.Range('E18').Value = 11509284
$ws.Range('F18').Value = 19752340
$ws.Range('G18').Value = 42493650
$ws.Range('H18').Value = 60361564
$ws.Range('D19').Value = 1155
$ws.Range('E19').Value = 160
$ws.Range('F19').Value = 25
$ws.Range('G19').Value = 0
$ws.Range('E20').Value = 365313
$ws.Range('F20').Value = 535130
$ws.Range('G20').Value = 365313
$ws.Range('H20').Value = 10703626
$ws.Range('D22').Value = 1381633
$ws.Range('E22').Value = 1647180
$ws.Range('F22').Value = 2368462
$ws.Range('G22').Value = 2741426
$ws.Range('H22').Value = 3365149
$ws.Range('D23').Value = 9039
$ws.Range('E23').Value = 8412
$ws.Range('F23').Value = 10220
$ws.Range('G23').Value = 12227
$ws.Range('H23').Value = 13091
$ws.Range('E25').Value = 11039
$ws.Range('H25').Value = 12039
$ws.Range('D26').Value = 1400277
$ws.Range('E26').Value = 2032104
$ws.Range('F26').Value = 2924876
$ws.Range('G26').Value = 3130005
$ws.Range('H26').Value = 14093905
$ws.Range('D27').Value = 8755655
$ws.Range('E27').Value = 13541388
$ws.Range('F27').Value = 22677216
$ws.Range('G27').Value = 45623655
$ws.Range('H27').Value = 74455469
$ws.Range('D29').Value = 113785
$ws.Range('E29').Value = 400939
$ws.Range('F29').Value = 1749781
$ws.Range('G29').Value = 626517
$ws.Range('H29').Value = 759304
$ws.Range('D31').Value = 320654
$ws.Range('E31').Value = 555877
$ws.Range('F31').Value = 871943
$ws.Range('G31').Value = 1423222
$ws.Range('H31').Value = 2064077
$ws.Range('D32').Value = 165936
$ws.Range('E32').Value = 332498
$ws.Range('F32').Value = 754977
$ws.Range('G32').Value = 2071884
$ws.Range('H32').Value = 6103818
$ws.Range('D33').Value = 18008
$ws.Range('E33').Value = 26428
$ws.Range('F33').Value = 31072
$ws.Range('G33').Value = 38300
$ws.Range('H33').Value = 110095
$ws.Range('D37').Value = 618383
$ws.Range('E37').Value = 1315742
$ws.Range('F37').Value = 3407773
$ws.Range('G37').Value = 4159923
$ws.Range('H37').Value = 9037294
$ws.Range('D39').Value = '-'
$ws.Range('D41').Value = 500183
$ws.Range('E41').Value = 671195
$ws.Range('F41').Value = 918334
$ws.Range('G41').Value = 1230418
$ws.Range('H41').Value = 1397869
$ws.Range('D42').Value = 500183
$ws.Range('E42').Value = 671195
$ws.Range('F42').Value = 918334
$ws.Range('G42').Value = 1230418
$ws.Range('H42').Value = 1397869
$ws.Range('D43').Value = 1118566
$ws.Range('E43').Value = 1986937
$ws.Range('F43').Value = 4326107
$ws.Range('G43').Value = 5390341
$ws.Range('H43').Value = 10435163
$ws.Range('F48').Value = -52126
$ws.Range('G48').Value = -76584
$ws.Range('H48').Value = -82828
$ws.Range('D49').Value = 0
$ws.Range('G49').Value = 11632
$ws.Range('H49').Value = 19590
$ws.Range('D52').Value = '-'
$ws.Range('D54').Value = '-'
$ws.Range('D56').Value = 5668186
$ws.Range('E56').Value = 9585548
$ws.Range('F56').Value = 16434332
$ws.Range('G56').Value = 38329363
$ws.Range('H56').Value = 62114641
$ws.Range('D57').Value = 7637089
$ws.Range('E57').Value = 11554451
$ws.Range('F57').Value = 18351109
$ws.Range('G57').Value = 40233314
$ws.Range('H57').Value = 64020306
$ws.Range('D58').Value = 8755655
$ws.Range('E58').Value = 13541388
$ws.Range('F58').Value = 22677216
$ws.Range('G58').Value = 45623655
$ws.Range('H58').Value = 74455469
